# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (column C) date value for every data row (2..31)
#    from 46070 to 46072 (serial date number).
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}

# 2) Rows 24-26 got cyclically re-ordered (row25 -> row24, row26 -> row25,
#    row24 -> row26), keeping columns A, B and G values together.
$row24A = $ws.Cells.Item(24, 1).Value2
$row24B = $ws.Cells.Item(24, 2).Value2
$row24G = $ws.Cells.Item(24, 7).Value2

$row25A = $ws.Cells.Item(25, 1).Value2
$row25B = $ws.Cells.Item(25, 2).Value2
$row25G = $ws.Cells.Item(25, 7).Value2

$row26A = $ws.Cells.Item(26, 1).Value2
$row26B = $ws.Cells.Item(26, 2).Value2
$row26G = $ws.Cells.Item(26, 7).Value2

# New row 24 <= old row 25
$ws.Cells.Item(24, 1).Value = $row25A
$ws.Cells.Item(24, 2).Value = $row25B
$ws.Cells.Item(24, 7).Value = $row25G

# New row 25 <= old row 26
$ws.Cells.Item(25, 1).Value = $row26A
$ws.Cells.Item(25, 2).Value = $row26B
$ws.Cells.Item(25, 7).Value = $row26G

# New row 26 <= old row 24
$ws.Cells.Item(26, 1).Value = $row24A
$ws.Cells.Item(26, 2).Value = $row24B
$ws.Cells.Item(26, 7).Value = $row24G
